$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.449.88'
$ws.Range("E2").Value = '  -0.07%  '

$ws.Range("D3").Value = '1.919.93'
$ws.Range("E3").Value = '  +0.83%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.010'
$ws.Range("E4").Value = '  +0.62%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.46'
$ws.Range("E5").Value = '  -0.06%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.008'
$ws.Range("E6").Value = '  +0.52%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4826'
$ws.Range("E7").Value = '  +0.57%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4069'
$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08215'
$ws.Range("E9").Value = '  +1.76%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.014'
$ws.Range("E10").Value = '  +1.04%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.27'
$ws.Range("E11").Value = '  -0.47%  '

$ws.Range("D12").Value = '1.907.05'
$ws.Range("E12").Value = '  -0.55%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.072'
$ws.Range("E13").Value = '  +1.95%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.245'
$ws.Range("E14").Value = '  +2.33%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.71'
$ws.Range("E15").Value = '  +1.78%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06872'
$ws.Range("E16").Value = '  +2.80%  '

$ws.Range("E17").Value = '  +0.50%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001040'
$ws.Range("E18").Value = '  +0.65%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.62'
$ws.Range("E19").Value = '  -0.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.009'
$ws.Range("E20").Value = '  +0.63%  '

$ws.Range("D21").Value = '29.448.53'
$ws.Range("E21").Value = '  -0.14%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.672'
$ws.Range("E22").Value = '  +2.33%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.75'
$ws.Range("E23").Value = '  -0.28%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.188'
$ws.Range("E24").Value = '  +1.06%  '

$ws.Range("D25").Value = '2.172.01'
$ws.Range("E25").Value = '  +1.39%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.647'
$ws.Range("E26").Value = '  +8.86%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.11'
$ws.Range("E27").Value = '  +1.01%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.02'
$ws.Range("E28").Value = '  +0.92%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.118'
$ws.Range("E29").Value = '  +1.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.79'
$ws.Range("E30").Value = '  +2.06%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.016'
$ws.Range("E31").Value = '  -1.73%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09633'
$ws.Range("E32").Value = '  +1.28%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.650'
$ws.Range("E33").Value = '  +2.40%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.549'
$ws.Range("E34").Value = '  +0.13%  '

$ws.Range("E35").Value = '  -1.48%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02285'
$ws.Range("E36").Value = '  +1.42%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06110'
$ws.Range("E37").Value = '  +0.55%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.185'
$ws.Range("E38").Value = '  +0.70%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '10.90'
$ws.Range("E39").Value = '  +6.67%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.068'
$ws.Range("E40").Value = '  +1.64%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5969'
$ws.Range("E41").Value = '  +1.10%  '

$ws.Range("E42").Value = '  +0.11%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.282'
$ws.Range("E43").Value = '  -0.18%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.392'
$ws.Range("E44").Value = '  -0.31%  '

$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.07613'
$ws.Range("E45").Value = '  -2.31%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.42'
$ws.Range("E46").Value = '  +0.87%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5595'
$ws.Range("E47").Value = '  +1.05%  '

$ws.Range("E48").Value = '  +1.49%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '118.91'
$ws.Range("E49").Value = '  +4.22%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.432'
$ws.Range("E50").Value = '  +3.93%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.29'
$ws.Range("E51").Value = '  -0.20%  '
